$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "HEX: #ffffff"
$ws.Range("B2").Value = "HEX: #ffffff"
$ws.Range("C2").Value = '<a class="color-link" href="/?themecolor=default">.</a>'
$ws.Range("D2").Value = "1:1"
$ws.Range("E2").Value = "Increase contrast by at least 350% to pass."

$ws.Range("A3").Value = "HEX: #000000"
$ws.Range("B3").Value = "HEX: #000000"
$ws.Range("C3").Value = '<a class="color-link" href="/?themecolor=blacktheme">.</a>'
$ws.Range("D3").Value = "1:1"
$ws.Range("E3").Value = "Increase contrast by at least 350% to pass."
